$d = $word.ActiveDocument

$titleText = "Play Coin-O-Mania for Free: Review of IGT's Treasure Hunt Adventure Slot"
$oldSummaryText = "Join a crew on a treasure hunt in Coin-O-Mania, the 5-reel, 4-row slot game from IGT. Play for free and read our review for more information."

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    Heading1 title paragraph at the top of the document:
#      <empty run><bold "Meta description" run><plain ": ..." run>
# -----------------------------------------------------------------
$titlePara = $null
$n0 = $d.Paragraphs.Count
for ($i = 1; $i -le $n0; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd() -eq $titleText) {
        $titlePara = $cand
        break
    }
}

$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $titlePara.Next()
$metaPara.Style = "Normal"
$metaRange = $metaPara.Range

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Join a crew on a treasure hunt in Coin-O-Mania, the 5-reel, 4-row slot game from IGT. Play for free and read our review for more information.</w:t></w:r></w:p>'
$metaRange.InsertXML($metaXml)

# -----------------------------------------------------------------
# 2) Near the end of the document, drop the duplicated bold title
#    paragraph ("Play Coin-O-Mania for Free: ...") that used to sit
#    just before the italic summary paragraph (search from the end
#    so we find the *other* copy of the title, not the one at the
#    top of the document).
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd() -eq $titleText) {
        $cand.Range.Delete()
        break
    }
}

# -----------------------------------------------------------------
# 3) Replace the text of the trailing italic paragraph with the new
#    image-prompt copy, preserving its run/formatting structure and
#    avoiding any smart-quote auto-correction of the apostrophe.
# -----------------------------------------------------------------
$summaryPara = $null
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd() -eq $oldSummaryText) {
        $summaryPara = $cand
        break
    }
}

$summaryRange = $summaryPara.Range
$startPos = $summaryRange.Start
$oldEndPos = $summaryRange.End - 1

$newText = "Create a feature image for Coin-o-Mania that captures the playful and adventurous spirit of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses, as this is one of the game's unique and memorable symbols. The warrior should be surrounded by symbols from the game, such as a pirate ship, treasure chest, and golden coins. The background should be a vibrant and colorful ocean setting, with waves crashing and seagulls flying overhead. The overall effect should be lively, engaging, and inviting, enticing players to join the adventure and discover their own treasures on the high seas."

# Insert the new text (plain, formatting applied explicitly afterwards)
# right after the existing text, then delete the old text. Using
# InsertAfter (rather than Find/Replace) keeps straight apostrophes
# from being auto-corrected into curly/smart quotes.
$insPoint = $d.Range($oldEndPos, $oldEndPos)
$insPoint.InsertAfter($newText)

$newLen = $newText.Length
$newRange = $d.Range($oldEndPos, $oldEndPos + $newLen)
$newRange.Font.Italic = $true

$oldRange = $d.Range($startPos, $oldEndPos)
$oldRange.Delete()

Write-Output "done"
